$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old first data row (A2:C2); remaining data rows shift up by one.
$ws.Rows(2).Delete()

# Append the newly captured measurements (11 rows) after the now-last row (21).
$ws.Range("A21").Value = -0.2765692472457886
$ws.Range("B21").Value = 2.626567840576172
$ws.Range("C21").Value = 0.6087272167205811

$ws.Range("A22").Value = -0.3174972236156463
$ws.Range("B22").Value = 1.27487576007843
$ws.Range("C22").Value = 0.1149953827261924

$ws.Range("A23").Value = -0.2593123018741607
$ws.Range("B23").Value = -0.741895854473114
$ws.Range("C23").Value = -0.290313720703125

$ws.Range("A24").Value = 0.6250678896903992
$ws.Range("B24").Value = -3.214983701705933
$ws.Range("C24").Value = -0.8246681094169617

$ws.Range("A25").Value = 0.2973386645317077
$ws.Range("B25").Value = -4.603633403778076
$ws.Range("C25").Value = 0.0426078513264656

$ws.Range("A26").Value = -0.3110831379890442
$ws.Range("B26").Value = -4.220009803771973
$ws.Range("C26").Value = 1.416138410568237

$ws.Range("A27").Value = -0.2417499274015426
$ws.Range("B27").Value = -3.320205211639404
$ws.Range("C27").Value = 1.446528911590576

$ws.Range("A28").Value = 0.3194825351238251
$ws.Range("B28").Value = -0.44180828332901
$ws.Range("C28").Value = 0.5250386595726013

$ws.Range("A29").Value = 0.0740674138069152
$ws.Range("B29").Value = 2.8290696144104
$ws.Range("C29").Value = -0.7684684991836548

$ws.Range("A30").Value = 0.3081815242767334
$ws.Range("B30").Value = 5.233893394470215
$ws.Range("C30").Value = -0.9094256162643432

$ws.Range("A31").Value = -0.4476115107536316
$ws.Range("B31").Value = 4.42144250869751
$ws.Range("C31").Value = 0.3985895812511444
